# Amend data to include ID column
# Adds an "ID" column (F) to Sheet 3, numbering each of the two 30-row
# plate blocks 1-30, and makes Sheet 3 the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Header cell
$ws3.Cells.Item(1, 6).Value = "ID"

# First plate block: rows 2-31 -> ID 1-30
for ($i = 1; $i -le 30; $i++) {
    $ws3.Cells.Item($i + 1, 6).Value = $i
}

# Second plate block: rows 32-61 -> ID 1-30
for ($i = 1; $i -le 30; $i++) {
    $ws3.Cells.Item($i + 31, 6).Value = $i
}

# Make Sheet 3 the active sheet/tab and set its selection + scroll position
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("H33").Select()
